# Atualização de bases das ligas, do dia: 14-04-2024 às 18:28
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Rows 11 and 12: the two fixtures had their data swapped (ids in
# column A stay at 9 / 10, everything else from B..AC is exchanged
# between the rows).
# ---------------------------------------------------------------

# New values for row 11 (previously held by row 12)
$ws.Range("B11").Value2 = 6240280
$ws.Range("F11").Value2 = "Atletico Ottawa"
$ws.Range("G11").Value2 = "Vancouver FC"
$ws.Range("H11").Value2 = 1
$ws.Range("I11").Value2 = 0
$ws.Range("K11").Value2 = 1.571
$ws.Range("L11").Value2 = 3.4
$ws.Range("M11").Value2 = 5.5
$ws.Range("N11").Value2 = 1.444
$ws.Range("O11").Value2 = 3.8
$ws.Range("P11").Value2 = 6
$ws.Range("Q11").Value2 = -1.25
$ws.Range("R11").Value2 = 1.95
$ws.Range("S11").Value2 = 1.85
$ws.Range("T11").Value2 = 2.75
$ws.Range("U11").Value2 = 1.975
$ws.Range("V11").Value2 = 1.825
$ws.Range("W11").Value2 = 0.444
$ws.Range("Z11").Value2 = -0.5
$ws.Range("AA11").Value2 = 0.425
$ws.Range("AB11").Value2 = -1
$ws.Range("AC11").Value2 = 0.825

# New values for row 12 (previously held by row 11)
$ws.Range("B12").Value2 = 6227815
$ws.Range("F12").Value2 = "HFX Wanderers"
$ws.Range("G12").Value2 = "Cavalry FC"
$ws.Range("H12").Value2 = 3
$ws.Range("I12").Value2 = 1
$ws.Range("K12").Value2 = 2.6
$ws.Range("L12").Value2 = 3.2
$ws.Range("M12").Value2 = 2.4
$ws.Range("N12").Value2 = 3.3
$ws.Range("O12").Value2 = 3
$ws.Range("P12").Value2 = 2.15
$ws.Range("Q12").Value2 = 0.25
$ws.Range("R12").Value2 = 1.925
$ws.Range("S12").Value2 = 1.875
$ws.Range("T12").Value2 = 2.25
$ws.Range("U12").Value2 = 2
$ws.Range("V12").Value2 = 1.8
$ws.Range("W12").Value2 = 2.3
$ws.Range("Z12").Value2 = 0.925
$ws.Range("AA12").Value2 = -1
$ws.Range("AB12").Value2 = 1
$ws.Range("AC12").Value2 = -1

# ---------------------------------------------------------------
# Row 89: fixture id / kickoff date / teams and odds refreshed
# ---------------------------------------------------------------
$ws.Range("B89").Value2 = 7803362
$ws.Range("E89").Value2 = 45396.83333333334
$ws.Range("F89").Value2 = "Vancouver FC"
$ws.Range("G89").Value2 = "Valour FC"
$ws.Range("K89").Value2 = 2.4
$ws.Range("M89").Value2 = 2.4
$ws.Range("N89").Value2 = 3
$ws.Range("P89").Value2 = 2
$ws.Range("Q89").Value2 = 0.25
$ws.Range("R89").Value2 = 2
$ws.Range("S89").Value2 = 1.8
$ws.Range("U89").Value2 = 1.875
$ws.Range("V89").Value2 = 1.925
